# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest count) figures and a refreshed cover
# image URL to the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 52
$ws1.Range("F4").Value = 2293
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg"
$ws1.Range("F5").Value = 22
$ws1.Range("F6").Value = 504

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 52
$ws4.Range("F6").Value = 2293
$ws4.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg"
$ws4.Range("F7").Value = 22
$ws4.Range("F8").Value = 504
